# fixbug tinh chiet khau don thu no
# Update static salary summary values on the per-employee sheets affected by
# the debt-collection-discount calculation fix.

$wb = $excel.ActiveWorkbook

# NV-11 Đỗ Thị Huyền Trân
$ws = $wb.Worksheets.Item("NV-11 Đỗ Thị Huyền Trân")
$ws.Range("B11").Value = 80000
$ws.Range("B35").Value = 9327714.285714285
$ws.Range("B38").Value = 9327714.285714285

# NV-29 Lâm Hoàng Phú
$ws = $wb.Worksheets.Item("NV-29 Lâm Hoàng Phú")
$ws.Range("B10").Value = 170000
$ws.Range("B34").Value = 823571.4285714286
$ws.Range("B37").Value = 923571.4285714286

# NV-5 Nguyễn Hoàng Yến Quyên
$ws = $wb.Worksheets.Item("NV-5 Nguyễn Hoàng Yến Quyên")
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 1142857.142857143
$ws.Range("B17").Value = 3350000
$ws.Range("B22").Value = 5
$ws.Range("B23").Value = 1428571.428571429
$ws.Range("B33").Value = 4492857.142857143
$ws.Range("B34").Value = 5248571.428571429
$ws.Range("B35").Value = 23745285.71428572

# NV-7 Phạm Thanh Hoàng
$ws = $wb.Worksheets.Item("NV-7 Phạm Thanh Hoàng")
$ws.Range("B25").Value = 1
$ws.Range("B26").Value = 321428.5714285714
$ws.Range("B38").Value = 5032142.857142857
$ws.Range("B39").Value = 14027500

# NV-22 Nguyễn Phúc Nam
$ws = $wb.Worksheets.Item("NV-22 Nguyễn Phúc Nam")
$ws.Range("B12").Value = 24.5
$ws.Range("B13").Value = 7000000.000000001
$ws.Range("B14").Value = 2600000
$ws.Range("B32").Value = 600000
$ws.Range("B34").Value = 600000

# NV-23 Lê Hoàng Thanh
$ws = $wb.Worksheets.Item("NV-23 Lê Hoàng Thanh")
$ws.Range("B12").Value = 25
$ws.Range("B13").Value = 4464285.714285715
$ws.Range("B32").Value = 3464285.714285715
$ws.Range("B34").Value = 4164285.714285715

# NV-30 Đào Vương Anh
$ws = $wb.Worksheets.Item("NV-30 Đào Vương Anh")
$ws.Range("B12").Value = 23
$ws.Range("B13").Value = 3285714.285714286
$ws.Range("B18").Value = 700000
$ws.Range("B32").Value = 985714.2857142859
$ws.Range("B34").Value = 985714.2857142859

# NV-36 Đặng Ngọc Mai
$ws = $wb.Worksheets.Item("NV-36 Đặng Ngọc Mai")
$ws.Range("B13").Value = 23
$ws.Range("B14").Value = 805000
$ws.Range("B15").Value = 4107142.857142857
$ws.Range("B36").Value = 6112142.857142857
$ws.Range("B38").Value = 6112142.857142857

# NV-40 Sang sang
$ws = $wb.Worksheets.Item("NV-40 Sang sang")
$ws.Range("B12").Value = 21
$ws.Range("B13").Value = 735000
$ws.Range("B14").Value = 2250000
$ws.Range("B33").Value = 1655000
$ws.Range("B35").Value = 1655000
